$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits -------------------------------------------------------
# Browser changed from Edge to Chrome
$ws.Range("D2").Value = "Chrome"

# AddFile / FileUpload rows: Runmode flipped from Yes to No
$ws.Range("F17").Value = "No"
$ws.Range("F18").Value = "No"

# Quit -> Close (and its Runmode flipped from Yes to No)
$ws.Range("C28").Value = "Close"
$ws.Range("F28").Value = "No"

# New rows 22 & 23 (logout flow): fill E22 before C22 so the shared-string
# table gets built in the same order as the authored workbook.
$ws.Range("E22").Value = "TryLogout"
$ws.Range("C22").Value = "MouseClick"
$ws.Range("D22").Value = "Nil"
$ws.Range("F22").Value = "Yes"

# C23 reuses the plain bordered style already used throughout the sheet
# (same as C22) — copy that formatting across before writing the value.
$ws.Range("C22").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = "MouseClick"
$ws.Range("D23").Value = "Nil"
$ws.Range("E23").Value = "Logout"
$ws.Range("F23").Value = "Yes"

# --- Formatting: give D23, E23 & F23 each a left/right thin border -------
# (no top/bottom) — apply per-cell so every cell gets both edges, not just
# the outer boundary of the combined range.
foreach ($addr in @("D23", "E23", "F23")) {
    $cell = $ws.Range($addr)
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
}

# --- Data validation: dropdown list on D2 sourced from $H$2:$H$4 ---------
$ws.Range("D2").Validation.Add(3, 1, 1, "=`$H`$2:`$H`$4")

# --- Selection state -------------------------------------------------------
$ws.Range("D23").Select()

Write-Output "done"
